$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "'2020.05.05"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").NumberFormat = "h:mm"
$ws.Range("B10").Value = 0.63888888888888895
$ws.Range("C10").NumberFormat = "h:mm"
$ws.Range("C10").Value = 0.67013888888888884
$ws.Range("E10").WrapText = $true
$ws.Range("E10").Value = "recreate phasor and edge, which are MSP objects (as opposed to Max objects)"

# Row 11
$ws.Rows.Item(11).RowHeight = 30
$ws.Range("B11").NumberFormat = "h:mm"
$ws.Range("B11").Value = 0.68055555555555547
$ws.Range("C11").NumberFormat = "h:mm"
$ws.Range("C11").Value = 0.77430555555555547
$ws.Range("E11").WrapText = $true
$ws.Range("E11").Value = "now we delve into buffers. Exploring buffer~ object and index~ object first. Lots of useful guides in source >> min_api >> doc"

# Row 12
$ws.Range("B12").NumberFormat = "h:mm"
$ws.Range("B12").Value = 0.875
$ws.Range("E12").WrapText = $true
$ws.Range("E12").Value = "merge branch, get latest Resonance, create new branch and start building something that includes Res"

# Update the active selection to match the edit
$ws.Range("E13").Select()
